$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
# fixed_ticker(I) is set to "LPSN" for every data row (restores the correct ticker after a bad join).
$rowsData = @(
    @(2, 152.8500061035156, 141, 156.1499938964844, 134.8500061035156, 10039606),
    @(3, 148.3500061035156, 144.1499938964844, 154.0500030517578, 117.75, 10039606),
    @(4, 112.1999969482422, 117, 124.0500030517578, 108.3000030517578, 10039606),
    @(5, 100.9499969482422, 84.90000152587891, 102.9000015258789, 82.94999694824219, 10039606),
    @(6, 87.30000305175781, 90.75, 98.5500030517578, 84.15000152587891, 10039606),
    @(7, 95.40000152587891, 100.1999969482422, 117.4499969482422, 92.6999969482422, 10039606),
    @(8, 125.6999969482422, 127.5, 131.1000061035156, 122.8499984741211, 10039606),
    @(9, 111.75, 109.5, 123, 102.75, 10039606),
    @(10, 101.25, 105.75, 110.25, 97.1999969482422, 10039606),
    @(11, 165, 203.25, 210.75, 160.5, 10039606),
    @(12, 203.25, 210.75, 226.5, 203.25, 10039606),
    @(13, 174, 179.25, 191.25, 171.75, 10039606),
    @(14, 245.25, 252, 257.25, 222.1499938964844, 10039606),
    @(15, 315, 348, 370.9500122070313, 312.75, 10039606),
    @(16, 392.8500061035156, 339, 393.6000061035156, 303.6000061035156, 10039606),
    @(17, 277.5, 352.0499877929688, 354.1499938964844, 266.7000122070312, 10039606),
    @(18, 440.8500061035156, 439.9500122070313, 453.2999877929688, 417.8999938964844, 10039606),
    @(19, 428.25, 497.8500061035156, 501.6000061035156, 414.2999877929688, 10039606),
    @(20, 535.5, 615.75, 622.7999877929688, 511.2000122070313, 10039606),
    @(21, 562.2000122070312, 615.1500244140625, 660.2999877929688, 540.1500244140625, 10039606),
    @(22, 334.2000122070312, 359.1000061035156, 378.4500122070313, 295.9500122070312, 10039606),
    @(23, 624.5999755859375, 644.7000122070312, 680.25, 556.5, 10039606),
    @(24, 803.4000244140625, 801.9000244140625, 930.2999877929688, 766.2000122070312, 10039606),
    @(25, 936.9000244140624, 950.4000244140624, 1032.449951171875, 889.7999877929688, 10039606),
    @(26, 788.25, 819.75, 893.8499755859375, 788.25, 10039606),
    @(27, 933.4500122070312, 955.3499755859376, 969.1500244140624, 814.7999877929688, 10039606),
    @(28, 887.8499755859375, 772.6500244140625, 925.5, 765.4500122070312, 10039606),
    @(29, 537.5999755859375, 448.0499877929688, 539.4000244140625, 377.8500061035156, 10039606),
    @(30, 366.8999938964844, 339.2999877929688, 399.8999938964844, 315.2999877929688, 10039606),
    @(31, 213.75, 204.6000061035156, 250.5, 188.1000061035156, 10039606),
    @(32, 143.8500061035156, 158.5500030517578, 163.1999969482422, 119.4000015258789, 10039606),
    @(33, 157.8000030517578, 193.1999969482422, 196.5, 145.6499938964844, 10039606),
    @(34, 66.44999694824219, 69.44999694824219, 84.30000305175781, 61.65000152587891, 10039606),
    @(35, 67.05000305175781, 71.25, 96.15000152587891, 60.45000076293945, 10039606),
    @(36, 59.09999847412109, 39.59999847412109, 59.40000152587891, 34.95000076293945, 10039606),
    @(37, 56.84999847412109, 42, 59.70000076293945, 41.40000152587891, 10039606),
    @(38, 15.14999961853027, 7.5, 15.14999961853027, 6.75, 10039606),
    @(39, 8.850000381469727, 18.14999961853028, 24, 8.550000190734863, 10039606),
    @(40, 19.20000076293945, 17.70000076293945, 19.35000038146973, 15.89999961853027, 10039606),
    @(41, 23.25, 21.60000038146973, 28.20000076293945, 18.14999961853028, 10039606),
    @(42, 12, 13.05000019073486, 13.5, 9.149999618530272, 10039606),
    @(43, 14.85000038146973, 14.10000038146973, 17.70000076293945, 13.19999980926514, 10039606)
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Range("D$r").Value = $row[1]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = $row[5]
    $ws.Range("I$r").Value = "LPSN"
}